$d = $word.ActiveDocument

$replacements = @(
    @{old="126÷7="; new="515÷7="},
    @{old="790÷4="; new="424÷2="},
    @{old="753÷7="; new="837÷3="},
    @{old="869÷7="; new="823÷4="},
    @{old="113÷8="; new="530÷7="},
    @{old="149÷3="; new="153÷7="},
    @{old="209÷7="; new="729÷8="},
    @{old="552÷4="; new="856÷8="},
    @{old="264÷9="; new="319÷7="},
    @{old="712÷4="; new="806÷9="},
    @{old="903÷5="; new="962÷8="},
    @{old="329÷2="; new="393÷8="},
    @{old="503÷4="; new="761÷4="},
    @{old="492÷6="; new="307÷8="},
    @{old="632÷9="; new="511÷8="},
    @{old="938÷5="; new="497÷5="},
    @{old="386÷8="; new="845÷3="},
    @{old="660÷4="; new="661÷5="},
    @{old="974÷7="; new="232÷2="},
    @{old="459÷4="; new="280÷2="},
    @{old="675÷3="; new="715÷4="},
    @{old="500÷7="; new="132÷3="},
    @{old="954÷9="; new="233÷6="},
    @{old="996÷9="; new="995÷4="},
    @{old="392÷9="; new="408÷8="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
